$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-19 Friday" "2024-04-20 Saturday"

Replace-Text "138÷6=" "789÷9="
Replace-Text "691÷3=" "946÷8="
Replace-Text "523÷4=" "316÷4="
Replace-Text "216÷3=" "416÷8="
Replace-Text "778÷4=" "719÷9="

Replace-Text "463÷2=" "203÷7="
Replace-Text "954÷9=" "855÷9="
Replace-Text "457÷4=" "657÷9="
Replace-Text "999÷4=" "479÷6="
Replace-Text "192÷2=" "858÷2="

Replace-Text "286÷7=" "627÷6="
Replace-Text "461÷6=" "540÷7="
Replace-Text "630÷2=" "445÷9="
Replace-Text "128÷6=" "173÷4="
Replace-Text "287÷6=" "672÷4="

Replace-Text "777÷5=" "108÷9="
Replace-Text "448÷7=" "210÷5="
Replace-Text "973÷7=" "511÷7="
Replace-Text "141÷8=" "550÷5="
Replace-Text "803÷4=" "997÷9="

Replace-Text "877÷4=" "455÷9="
Replace-Text "569÷6=" "720÷4="
Replace-Text "397÷5=" "886÷4="
Replace-Text "315÷5=" "220÷3="
Replace-Text "970÷4=" "534÷7="
